$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("30.09.")
$ws3 = $wb.Worksheets.Item("Notizen Präsi")

# New shared-string values are entered in this order so the shared string
# table ends up with the same ordering as the target workbook:
# 49 await war wegen promise?         (sheet3 A4)
# 50 filter in Liste und Archiv funktionieren (sheet2 A2)
# 51 in Liste bearbeiten und delte button bearbeitet (sheet2 A3)
# 52 reactive forms                   (sheet3 A5)
# 53 activated route                  (sheet3 A6)
# 54 router?                          (sheet3 A7)
$ws3.Range("A4").Value = "await war wegen promise?"

# Sheet2 ("30.09.") - add two new rows of notes
$ws2.Range("A2").Value = "filter in Liste und Archiv funktionieren"
$ws2.Range("A3").Value = "in Liste bearbeiten und delte button bearbeitet"
$ws2.Range("A3").Select()

# Sheet3 ("Notizen Präsi") - add remaining rows of notes
$ws3.Range("A5").Value = "reactive forms"
$ws3.Range("A6").Value = "activated route"
$ws3.Range("A7").Value = "router?"

# Sheet1 ("Tabelle1") - move selection, no longer the active tab
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1.Range("E2").Select()

# Make "Notizen Präsi" the active sheet/tab, with selection on A8
$ws3.Activate()
$ws3.Range("A8").Select()

$wb.Save()
